# Daily attendance processing - 2026-01-27 09:18:03
# Normalises the "Recorded By" (column G) entries so that the
# first two comma-separated contributors are listed with the
# human/service account before the literal "System" token.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column G ("Recorded By") whose value needs the first two
# comma-separated entries swapped.
$rows = @(2, 4, 5, 8, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 26, 28, 30, 31, 34, 36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50, 52, 54, 56, 57, 60, 62, 63, 64, 65, 66, 67, 69, 70, 71, 72, 73, 74, 76, 78, 80, 81, 82, 83, 84, 85, 86, 87, 90, 92, 93, 94, 96, 99, 101, 106, 107, 108, 109, 110, 111, 112, 113, 116, 118, 119, 120, 122, 125, 127, 132, 133, 134, 135, 136, 137, 138, 139, 142, 144, 145, 146, 148, 151, 153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    $parts = $val.Split(",")

    $first = $parts[0].Trim()
    $second = $parts[1].Trim()

    if ($parts.Length -gt 2) {
        $third = $parts[2].Trim()
        $newVal = $second + ", " + $first + ", " + $third
    } else {
        $newVal = $second + ", " + $first
    }

    $cell.Value = $newVal
}
